# Insert a new data row for "La Paz, Bolivia" (LPB) immediately before the
# existing "Amman, Jordan" row (row 218), shifting rows 218:312 down to
# 219:313. This matches the target diff: dimension grows from A1:G312 to
# A1:G313 and every row from 218 onward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 218..312 down by one to make room for the new row.
$ws.Rows.Item(218).Insert()

# The freshly inserted row has no explicit style; copy the formatting
# (centered/bold/bordered colo-code style) from the row right below it
# (which used to be row 218, now row 219) onto the new row 218 so the
# "colo" column (A) keeps its original look (style index 1).
$ws.Range("A219").Copy()
$ws.Range("A218").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A218").Value = "LPB"
$ws.Range("B218").Value = "La Paz, Bolivia"
$ws.Range("C218").Value = -16.4897
$ws.Range("D218").Value = -68.1193
$ws.Range("E218").Value = "BO"
$ws.Range("F218").Value = "South America"
$ws.Range("G218").Value = "La Paz"
